# Apply updated crypto price/volume figures to Sheet1 (D = Price, E = Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "58.640.32"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.303.05"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.85"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.12"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  +2.87%  "
$ws.Range("D9").Value = "2.302.07"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.75"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "2.711.54"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "58.539.15"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "2.302.67"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "316.30"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.60"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.09"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.96"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.30"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.55"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").Value = "0.0₃0726"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.82"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.91"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "290.57"
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.11"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0952"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.557"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.29"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.63"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  +0.00%  "
